$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')
$ws.Cells.Item(1,1).Value = 0
$ws.Cells.Item(1,2).Value = '开始时间'
$ws.Cells.Item(1,3).Value = '名称'
$ws.Cells.Item(1,4).Value = '地点'
$ws.Cells.Item(1,5).Value = '具体时间范围'
$ws.Cells.Item(1,6).Value = '想去人数'
$ws.Cells.Item(1,7).Value = '最低票价'
$ws.Cells.Item(1,8).Value = 'Link'
$ws.Cells.Item(1,9).Value = 'Cover'
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = '2024-03-08'
$ws.Cells.Item(2,3).Value = '苏州·国风宠物-cosplay展（取消）'
$ws.Cells.Item(2,4).Value = '金山南路影视城 木渎影视城会展中心'
$ws.Cells.Item(2,5).Value = '2024.03.08 09:00-03.10 17:30'
$ws.Cells.Item(2,6).Value = 1166
$ws.Cells.Item(2,7).Value = '不可售'
$ws.Cells.Item(2,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80635'
$ws.Cells.Item(2,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg'
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = '2024-03-16'
$ws.Cells.Item(3,3).Value = '苏州·OrangeOrange新春随舞派对【免费展会】'
$ws.Cells.Item(3,4).Value = '狮山路298号 金鹰国际购物中心(狮山路店)'
$ws.Cells.Item(3,5).Value = '2024.03.16 13:00-03.16 17:30'
$ws.Cells.Item(3,6).Value = 74
$ws.Cells.Item(3,7).Value = 25
$ws.Cells.Item(3,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82046'
$ws.Cells.Item(3,9).Value = '//i2.hdslb.com/bfs/openplatform/202402/0OH3Ax4I1708913393518.png'
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = '2024-03-17'
$ws.Cells.Item(4,3).Value = '苏州·世纪幻想动漫游戏展2.0'
$ws.Cells.Item(4,4).Value = '清禾路886号 尹山湖大剧院'
$ws.Cells.Item(4,5).Value = '2024.03.17 10:00-03.17 17:00'
$ws.Cells.Item(4,6).Value = 1510
$ws.Cells.Item(4,7).Value = 60
$ws.Cells.Item(4,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81387'
$ws.Cells.Item(4,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/isVyI9hH1708590817616.jpeg'
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = '2024-03-23'
$ws.Cells.Item(5,3).Value = '苏州·Look Look动漫嘉年华'
$ws.Cells.Item(5,4).Value = '阳澄半岛慈云路168号(重元寺北) 阳澄湖澜廷度假酒店'
$ws.Cells.Item(5,5).Value = '2024.03.23 10:00-03.23 17:30'
$ws.Cells.Item(5,6).Value = 580
$ws.Cells.Item(5,7).Value = 52.2
$ws.Cells.Item(5,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81698'
$ws.Cells.Item(5,9).Value = '//i1.hdslb.com/bfs/openplatform/202402/CP95X8ao1708934930351.jpeg'
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = '2024-03-30'
$ws.Cells.Item(6,3).Value = '苏州·奇幻世界5.3动漫游戏展'
$ws.Cells.Item(6,4).Value = '龙河路1288号 乐动力苏州湾体育中心'
$ws.Cells.Item(6,5).Value = '2024.03.30 10:00-03.31 17:00'
$ws.Cells.Item(6,6).Value = 1069
$ws.Cells.Item(6,7).Value = 55
$ws.Cells.Item(6,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82002'
$ws.Cells.Item(6,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/HlxVHAz91708593664222.jpeg'
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = '2024-04-04'
$ws.Cells.Item(7,3).Value = '【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会'
$ws.Cells.Item(7,4).Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws.Cells.Item(7,5).Value = '2024.04.04 10:00-04.05 17:00'
$ws.Cells.Item(7,6).Value = 11080
$ws.Cells.Item(7,7).Value = 60
$ws.Cells.Item(7,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81827'
$ws.Cells.Item(7,9).Value = '//i2.hdslb.com/bfs/openplatform/202402/6oSFbWOx1707301464970.jpeg'
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = '2024-04-06'
$ws.Cells.Item(8,3).Value = '苏州·第一届寒假动漫展宅舞比赛-CF01'
$ws.Cells.Item(8,4).Value = '润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店'
$ws.Cells.Item(8,5).Value = '2024.04.06 10:00-04.06 16:00'
$ws.Cells.Item(8,6).Value = 85
$ws.Cells.Item(8,7).Value = 49
$ws.Cells.Item(8,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80528'
$ws.Cells.Item(8,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg'
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = '2024-04-13'
$ws.Cells.Item(9,3).Value = '苏州·X-party 国漫游戏嘉年华03'
$ws.Cells.Item(9,4).Value = '秋枫街与开平路交叉口西南角 爱琴海购物中心'
$ws.Cells.Item(9,5).Value = '2024.04.13 10:00-04.14 17:00'
$ws.Cells.Item(9,6).Value = 32
$ws.Cells.Item(9,7).Value = 48
$ws.Cells.Item(9,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82042'
$ws.Cells.Item(9,9).Value = '//i2.hdslb.com/bfs/openplatform/202403/GWNvc78z1709275224442.jpeg'
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = '2024-04-13'
$ws.Cells.Item(10,3).Value = '苏州·绘时国乙1.0-秩序之外'
$ws.Cells.Item(10,4).Value = '石路步行街永福桥浜15号 银河广场'
$ws.Cells.Item(10,5).Value = '2024.04.13 13:30-04.13 20:00'
$ws.Cells.Item(10,6).Value = 321
$ws.Cells.Item(10,7).Value = 88
$ws.Cells.Item(10,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80789'
$ws.Cells.Item(10,9).Value = '//i0.hdslb.com/bfs/openplatform/202403/nIPoXWqO1709275656198.jpeg'
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = '2024-04-20'
$ws.Cells.Item(11,3).Value = '苏州·首届Redamancy动漫游戏嘉年华'
$ws.Cells.Item(11,4).Value = '清禾路886号 尹山湖大剧院'
$ws.Cells.Item(11,5).Value = '2024.04.20 10:00-04.20 17:00'
$ws.Cells.Item(11,6).Value = 1067
$ws.Cells.Item(11,7).Value = 60
$ws.Cells.Item(11,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81879'
$ws.Cells.Item(11,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/lR4oJWzI1708309129629.jpeg'
$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = '2024-04-21'
$ws.Cells.Item(12,3).Value = '苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0'
$ws.Cells.Item(12,4).Value = '清禾路888号2号楼3楼 格莱美婚礼宴会中心'
$ws.Cells.Item(12,5).Value = '2024.04.21 10:00-04.21 21:00'
$ws.Cells.Item(12,6).Value = 759
$ws.Cells.Item(12,7).Value = 69.90000000000001
$ws.Cells.Item(12,8).Value = 'https://show.bilibili.com/platform/detail.html?id=78666'
$ws.Cells.Item(12,9).Value = '//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg'
$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = '2024-05-01'
$ws.Cells.Item(13,3).Value = '昆山·第十二届理想乡动漫游戏展'
$ws.Cells.Item(13,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(13,5).Value = '2024.05.01 10:00-05.03 17:00'
$ws.Cells.Item(13,6).Value = 12231
$ws.Cells.Item(13,7).Value = 75
$ws.Cells.Item(13,8).Value = 'https://show.bilibili.com/platform/detail.html?id=77196'
$ws.Cells.Item(13,9).Value = '//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png'
$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = '2024-05-01'
$ws.Cells.Item(14,3).Value = '苏州·第十七届 I COME ACG  动漫品牌博览会'
$ws.Cells.Item(14,4).Value = '金山南路288号 广电国际会展中心'
$ws.Cells.Item(14,5).Value = '2024.05.01 10:00-05.02 17:00'
$ws.Cells.Item(14,6).Value = 12775
$ws.Cells.Item(14,7).Value = 65
$ws.Cells.Item(14,8).Value = 'https://show.bilibili.com/platform/detail.html?id=79789'
$ws.Cells.Item(14,9).Value = '//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg'
$ws.Cells.Item(15,1).Value = 14
$ws.Cells.Item(15,2).Value = '2024-05-02'
$ws.Cells.Item(15,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会'
$ws.Cells.Item(15,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(15,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws.Cells.Item(15,6).Value = 32
$ws.Cells.Item(15,7).Value = 1
$ws.Cells.Item(15,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81116'
$ws.Cells.Item(15,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg'
$ws.Cells.Item(16,1).Value = 15
$ws.Cells.Item(16,2).Value = '2024-05-02'
$ws.Cells.Item(16,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会'
$ws.Cells.Item(16,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(16,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws.Cells.Item(16,6).Value = 124
$ws.Cells.Item(16,7).Value = 1
$ws.Cells.Item(16,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81100'
$ws.Cells.Item(16,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg'
$ws.Cells.Item(17,1).Value = 16
$ws.Cells.Item(17,2).Value = '2024-05-02'
$ws.Cells.Item(17,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会'
$ws.Cells.Item(17,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(17,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws.Cells.Item(17,6).Value = 18
$ws.Cells.Item(17,7).Value = 1
$ws.Cells.Item(17,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81119'
$ws.Cells.Item(17,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg'
$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = '2024-05-02'
$ws.Cells.Item(18,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会'
$ws.Cells.Item(18,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(18,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws.Cells.Item(18,6).Value = 29
$ws.Cells.Item(18,7).Value = 1
$ws.Cells.Item(18,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81118'
$ws.Cells.Item(18,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg'
$ws.Cells.Item(19,1).Value = 18
$ws.Cells.Item(19,2).Value = '2024-05-03'
$ws.Cells.Item(19,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会'
$ws.Cells.Item(19,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(19,5).Value = '2024.05.03 14:00-05.03 16:00'
$ws.Cells.Item(19,6).Value = 78
$ws.Cells.Item(19,7).Value = 1
$ws.Cells.Item(19,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81120'
$ws.Cells.Item(19,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg'
$ws.Cells.Item(20,1).Value = 19
$ws.Cells.Item(20,2).Value = '2024-05-03'
$ws.Cells.Item(20,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会'
$ws.Cells.Item(20,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(20,5).Value = '2024.05.03 14:00-05.03 16:00'
$ws.Cells.Item(20,6).Value = 43
$ws.Cells.Item(20,7).Value = 1
$ws.Cells.Item(20,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81114'
$ws.Cells.Item(20,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg'
$ws.Cells.Item(21,1).Value = 20
$ws.Cells.Item(21,2).Value = '2024-06-08'
$ws.Cells.Item(21,3).Value = '【会员购严选】苏州·Come in joy动漫国潮文化节'
$ws.Cells.Item(21,4).Value = '金山南路288号 广电国际会展中心'
$ws.Cells.Item(21,5).Value = '2024.06.08 10:00-06.09 17:00'
$ws.Cells.Item(21,6).Value = 43
$ws.Cells.Item(21,7).Value = 60
$ws.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82233'
$ws.Cells.Item(21,9).Value = '//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg'
$ws.Rows.Item(22).ClearContents()

$ws = $wb.Worksheets.Item('全部类型')
$ws.Cells.Item(1,1).Value = 0
$ws.Cells.Item(1,2).Value = '开始时间'
$ws.Cells.Item(1,3).Value = '名称'
$ws.Cells.Item(1,4).Value = '地点'
$ws.Cells.Item(1,5).Value = '具体时间范围'
$ws.Cells.Item(1,6).Value = '想去人数'
$ws.Cells.Item(1,7).Value = '最低票价'
$ws.Cells.Item(1,8).Value = 'Link'
$ws.Cells.Item(1,9).Value = 'Cover'
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = '2024-03-03'
$ws.Cells.Item(2,3).Value = '苏州·龙猫和他的朋友·动漫作品音乐会'
$ws.Cells.Item(2,4).Value = '星湖街555号高教区(体育馆南侧) 苏州独墅湖影剧院'
$ws.Cells.Item(2,5).Value = '2024.03.03 19:30-03.03 21:00'
$ws.Cells.Item(2,6).Value = 13
$ws.Cells.Item(2,7).Value = '不可售'
$ws.Cells.Item(2,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81799'
$ws.Cells.Item(2,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/gqnOEjvJ1707214629948.jpeg'
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = '2024-03-08'
$ws.Cells.Item(3,3).Value = '苏州·国风宠物-cosplay展（取消）'
$ws.Cells.Item(3,4).Value = '金山南路影视城 木渎影视城会展中心'
$ws.Cells.Item(3,5).Value = '2024.03.08 09:00-03.10 17:30'
$ws.Cells.Item(3,6).Value = 1166
$ws.Cells.Item(3,7).Value = '不可售'
$ws.Cells.Item(3,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80635'
$ws.Cells.Item(3,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg'
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = '2024-03-16'
$ws.Cells.Item(4,3).Value = '苏州·OrangeOrange新春随舞派对【免费展会】'
$ws.Cells.Item(4,4).Value = '狮山路298号 金鹰国际购物中心(狮山路店)'
$ws.Cells.Item(4,5).Value = '2024.03.16 13:00-03.16 17:30'
$ws.Cells.Item(4,6).Value = 74
$ws.Cells.Item(4,7).Value = 25
$ws.Cells.Item(4,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82046'
$ws.Cells.Item(4,9).Value = '//i2.hdslb.com/bfs/openplatform/202402/0OH3Ax4I1708913393518.png'
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = '2024-03-17'
$ws.Cells.Item(5,3).Value = '苏州·世纪幻想动漫游戏展2.0'
$ws.Cells.Item(5,4).Value = '清禾路886号 尹山湖大剧院'
$ws.Cells.Item(5,5).Value = '2024.03.17 10:00-03.17 17:00'
$ws.Cells.Item(5,6).Value = 1510
$ws.Cells.Item(5,7).Value = 60
$ws.Cells.Item(5,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81387'
$ws.Cells.Item(5,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/isVyI9hH1708590817616.jpeg'
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = '2024-03-23'
$ws.Cells.Item(6,3).Value = '苏州·Look Look动漫嘉年华'
$ws.Cells.Item(6,4).Value = '阳澄半岛慈云路168号(重元寺北) 阳澄湖澜廷度假酒店'
$ws.Cells.Item(6,5).Value = '2024.03.23 10:00-03.23 17:30'
$ws.Cells.Item(6,6).Value = 580
$ws.Cells.Item(6,7).Value = 52.2
$ws.Cells.Item(6,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81698'
$ws.Cells.Item(6,9).Value = '//i1.hdslb.com/bfs/openplatform/202402/CP95X8ao1708934930351.jpeg'
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = '2024-03-30'
$ws.Cells.Item(7,3).Value = '苏州·奇幻世界5.3动漫游戏展'
$ws.Cells.Item(7,4).Value = '龙河路1288号 乐动力苏州湾体育中心'
$ws.Cells.Item(7,5).Value = '2024.03.30 10:00-03.31 17:00'
$ws.Cells.Item(7,6).Value = 1069
$ws.Cells.Item(7,7).Value = 55
$ws.Cells.Item(7,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82002'
$ws.Cells.Item(7,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/HlxVHAz91708593664222.jpeg'
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = '2024-04-04'
$ws.Cells.Item(8,3).Value = '【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会'
$ws.Cells.Item(8,4).Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws.Cells.Item(8,5).Value = '2024.04.04 10:00-04.05 17:00'
$ws.Cells.Item(8,6).Value = 11080
$ws.Cells.Item(8,7).Value = 60
$ws.Cells.Item(8,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81827'
$ws.Cells.Item(8,9).Value = '//i2.hdslb.com/bfs/openplatform/202402/6oSFbWOx1707301464970.jpeg'
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = '2024-04-06'
$ws.Cells.Item(9,3).Value = '苏州·第一届寒假动漫展宅舞比赛-CF01'
$ws.Cells.Item(9,4).Value = '润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店'
$ws.Cells.Item(9,5).Value = '2024.04.06 10:00-04.06 16:00'
$ws.Cells.Item(9,6).Value = 85
$ws.Cells.Item(9,7).Value = 49
$ws.Cells.Item(9,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80528'
$ws.Cells.Item(9,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg'
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = '2024-04-13'
$ws.Cells.Item(10,3).Value = '苏州·X-party 国漫游戏嘉年华03'
$ws.Cells.Item(10,4).Value = '秋枫街与开平路交叉口西南角 爱琴海购物中心'
$ws.Cells.Item(10,5).Value = '2024.04.13 10:00-04.14 17:00'
$ws.Cells.Item(10,6).Value = 32
$ws.Cells.Item(10,7).Value = 48
$ws.Cells.Item(10,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82042'
$ws.Cells.Item(10,9).Value = '//i2.hdslb.com/bfs/openplatform/202403/GWNvc78z1709275224442.jpeg'
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = '2024-04-13'
$ws.Cells.Item(11,3).Value = '苏州·绘时国乙1.0-秩序之外'
$ws.Cells.Item(11,4).Value = '石路步行街永福桥浜15号 银河广场'
$ws.Cells.Item(11,5).Value = '2024.04.13 13:30-04.13 20:00'
$ws.Cells.Item(11,6).Value = 321
$ws.Cells.Item(11,7).Value = 88
$ws.Cells.Item(11,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80789'
$ws.Cells.Item(11,9).Value = '//i0.hdslb.com/bfs/openplatform/202403/nIPoXWqO1709275656198.jpeg'
$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = '2024-04-20'
$ws.Cells.Item(12,3).Value = '苏州·首届Redamancy动漫游戏嘉年华'
$ws.Cells.Item(12,4).Value = '清禾路886号 尹山湖大剧院'
$ws.Cells.Item(12,5).Value = '2024.04.20 10:00-04.20 17:00'
$ws.Cells.Item(12,6).Value = 1067
$ws.Cells.Item(12,7).Value = 60
$ws.Cells.Item(12,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81879'
$ws.Cells.Item(12,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/lR4oJWzI1708309129629.jpeg'
$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = '2024-04-21'
$ws.Cells.Item(13,3).Value = '苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0'
$ws.Cells.Item(13,4).Value = '清禾路888号2号楼3楼 格莱美婚礼宴会中心'
$ws.Cells.Item(13,5).Value = '2024.04.21 10:00-04.21 21:00'
$ws.Cells.Item(13,6).Value = 759
$ws.Cells.Item(13,7).Value = 69.90000000000001
$ws.Cells.Item(13,8).Value = 'https://show.bilibili.com/platform/detail.html?id=78666'
$ws.Cells.Item(13,9).Value = '//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg'
$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = '2024-05-01'
$ws.Cells.Item(14,3).Value = '昆山·第十二届理想乡动漫游戏展'
$ws.Cells.Item(14,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(14,5).Value = '2024.05.01 10:00-05.03 17:00'
$ws.Cells.Item(14,6).Value = 12231
$ws.Cells.Item(14,7).Value = 75
$ws.Cells.Item(14,8).Value = 'https://show.bilibili.com/platform/detail.html?id=77196'
$ws.Cells.Item(14,9).Value = '//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png'
$ws.Cells.Item(15,1).Value = 14
$ws.Cells.Item(15,2).Value = '2024-05-01'
$ws.Cells.Item(15,3).Value = '苏州·第十七届 I COME ACG  动漫品牌博览会'
$ws.Cells.Item(15,4).Value = '金山南路288号 广电国际会展中心'
$ws.Cells.Item(15,5).Value = '2024.05.01 10:00-05.02 17:00'
$ws.Cells.Item(15,6).Value = 12775
$ws.Cells.Item(15,7).Value = 65
$ws.Cells.Item(15,8).Value = 'https://show.bilibili.com/platform/detail.html?id=79789'
$ws.Cells.Item(15,9).Value = '//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg'
$ws.Cells.Item(16,1).Value = 15
$ws.Cells.Item(16,2).Value = '2024-05-02'
$ws.Cells.Item(16,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会'
$ws.Cells.Item(16,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(16,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws.Cells.Item(16,6).Value = 32
$ws.Cells.Item(16,7).Value = 1
$ws.Cells.Item(16,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81116'
$ws.Cells.Item(16,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg'
$ws.Cells.Item(17,1).Value = 16
$ws.Cells.Item(17,2).Value = '2024-05-02'
$ws.Cells.Item(17,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会'
$ws.Cells.Item(17,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(17,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws.Cells.Item(17,6).Value = 124
$ws.Cells.Item(17,7).Value = 1
$ws.Cells.Item(17,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81100'
$ws.Cells.Item(17,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg'
$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = '2024-05-02'
$ws.Cells.Item(18,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会'
$ws.Cells.Item(18,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(18,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws.Cells.Item(18,6).Value = 18
$ws.Cells.Item(18,7).Value = 1
$ws.Cells.Item(18,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81119'
$ws.Cells.Item(18,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg'
$ws.Cells.Item(19,1).Value = 18
$ws.Cells.Item(19,2).Value = '2024-05-02'
$ws.Cells.Item(19,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会'
$ws.Cells.Item(19,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(19,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws.Cells.Item(19,6).Value = 29
$ws.Cells.Item(19,7).Value = 1
$ws.Cells.Item(19,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81118'
$ws.Cells.Item(19,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg'
$ws.Cells.Item(20,1).Value = 19
$ws.Cells.Item(20,2).Value = '2024-05-03'
$ws.Cells.Item(20,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会'
$ws.Cells.Item(20,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(20,5).Value = '2024.05.03 14:00-05.03 16:00'
$ws.Cells.Item(20,6).Value = 78
$ws.Cells.Item(20,7).Value = 1
$ws.Cells.Item(20,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81120'
$ws.Cells.Item(20,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg'
$ws.Cells.Item(21,1).Value = 20
$ws.Cells.Item(21,2).Value = '2024-05-03'
$ws.Cells.Item(21,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会'
$ws.Cells.Item(21,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws.Cells.Item(21,5).Value = '2024.05.03 14:00-05.03 16:00'
$ws.Cells.Item(21,6).Value = 43
$ws.Cells.Item(21,7).Value = 1
$ws.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81114'
$ws.Cells.Item(21,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg'
$ws.Cells.Item(22,1).Value = 21
$ws.Cells.Item(22,2).Value = '2024-06-08'
$ws.Cells.Item(22,3).Value = '【会员购严选】苏州·Come in joy动漫国潮文化节'
$ws.Cells.Item(22,4).Value = '金山南路288号 广电国际会展中心'
$ws.Cells.Item(22,5).Value = '2024.06.08 10:00-06.09 17:00'
$ws.Cells.Item(22,6).Value = 43
$ws.Cells.Item(22,7).Value = 60
$ws.Cells.Item(22,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82233'
$ws.Cells.Item(22,9).Value = '//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg'
$ws.Rows.Item(23).ClearContents()

